# Fruta / hortaliza, semanal
# Rotate the data rows (2,3,4) down by one position, wrapping row 4 back to row 2:
#   old row 2 -> new row 3
#   old row 3 -> new row 4
#   old row 4 -> new row 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the full rows (columns A:R) before making any changes.
$row2 = $ws.Range("A2:R2").Value2
$row3 = $ws.Range("A3:R3").Value2
$row4 = $ws.Range("A4:R4").Value2

# Write them back in rotated order.
$ws.Range("A3:R3").Value2 = $row2
$ws.Range("A4:R4").Value2 = $row3
$ws.Range("A2:R2").Value2 = $row4
